$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 115, shifting existing rows 115:183 down to 116:183
$ws.Rows.Item(115).Insert()

# Populate the newly inserted row 115 with the new record
$ws.Range("A115").Value = 4
$ws.Range("B115").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C115").Value = "Los Lagos"
$ws.Range("D115").Value = 44452
$ws.Range("E115").Value = 10
$ws.Range("F115").Value = 100112008
$ws.Range("G115").Value = "Coliflor"
$ws.Range("H115").Value = "Sin especificar"
$ws.Range("I115").Value = "Segunda"
$ws.Range("J115").Value = 500
$ws.Range("K115").Value = 1000
$ws.Range("L115").Value = 1000
$ws.Range("M115").Value = 1000
$ws.Range("N115").Value = "$/unidad"
$ws.Range("O115").Value = "Región del Maule"
$ws.Range("P115").Value = 1000
$ws.Range("Q115").Value = 1
$ws.Range("R115").Value = "Hortaliza"
